# Auto-generated edit script applying the cryptos.xlsx diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '42.640.58'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -1.29%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.286.77'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -3.06%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.18%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '300.64'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -2.66%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '97.31'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -5.42%  '

$ws.Range("E7").Value = '  -1.42%  '

$ws.Range("E8").Value = '  -0.02%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.502'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -3.20%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '33.48'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -5.54%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0788'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -1.76%  '

$ws.Range("E12").Value = '  -4.54%  '

$ws.Range("E13").Value = '  -0.02%  '

$ws.Range("E14").Value = '  -3.15%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '2.641.35'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -3.29%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '15.36'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.78%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.298.87'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -2.70%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.788'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -2.16%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '42.548.60'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -1.48%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.0₃0896'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -1.55%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '11.49'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -2.57%  '

$ws.Range("E22").Value = '  -4.33%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '66.78'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -1.62%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '234.75'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -1.75%  '

$ws.Range("E25").Value = '  -4.16%  '

$ws.Range("E26").Value = '  -4.03%  '

$ws.Range("E27").Value = '  -0.03%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '24.47'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -4.11%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '166.12'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +2.98%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.07'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -6.54%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '33.79'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -6.25%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '9.11'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -3.16%  '

$ws.Range("E33").Value = '  -0.17%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.98'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -4.19%  '

$ws.Range("E35").Value = '  -3.81%  '

$ws.Range("E36").Value = '  -4.75%  '

$ws.Range("E37").Value = '  -6.35%  '

$ws.Range("B38").Value = 'LidoDAOToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.83'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -7.20%  '

$ws.Range("B39").Value = 'Celestia'
$ws.Range("C39").Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '16.25'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -9.65%  '

$ws.Range("E40").Value = '  -6.92%  '

$ws.Range("E41").Value = '  -3.80%  '

$ws.Range("E42").Value = '  -2.79%  '

$ws.Range("E43").Value = '  -4.39%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.963.37'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -3.59%  '

$ws.Range("E45").Value = '  -1.83%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '17.82'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -8.33%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '9.70'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -8.23%  '

$ws.Range("E48").Value = '  -7.24%  '

$ws.Range("B49").Value = 'HuobiToken'
$ws.Range("C49").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.82'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -3.01%  '

$ws.Range("B50").Value = 'MultiversX'
$ws.Range("C50").Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '53.20'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -7.35%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.507.48'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -3.37%  '
